$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Google Form entry-id mapping row with the new form's field IDs.
$ws.Range("B4").Value = "entry.98625466"
$ws.Range("C4").Value = "entry.1289724168"
$ws.Range("D4").Value = "entry.897012234"
$ws.Range("F4").Value = "entry.1846392603"
$ws.Range("H4").Value = "entry.1330307066"
$ws.Range("I4").Value = "entry.1030855690"
# Screenshot of Payment (image) field no longer mapped -> plain 0, like the other unused columns.
$ws.Range("M4").Value = 0

# Normalize the font back to the sheet default (drop the stray Consolas style).
$fmtSrc = $ws.Range("E4")
$fmtSrc.Copy()
$ws.Range("B4:D4").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("F4").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("H4:I4").PasteSpecial(-4122)
$fmtSrc.Copy()
$ws.Range("M4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

[void]$ws.Range("H10").Select()
